$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 157.89473
$ws.Range("I6").Value = 177.14285
$ws.Range("J6").Value = 104
$ws.Range("K6").Value = 531.4285500000001
$ws.Range("L6").Value = 312
$ws.Range("M6").Value = -419.4285500000001
$ws.Range("N6").Value = -536
$ws.Range("H8").Value = 3699
$ws.Range("I8").Value = 3699
$ws.Range("K8").Value = 11097
$ws.Range("M8").Value = -10958
$ws.Range("H17").Value = 498.27908
$ws.Range("J17").Value = 498.27908
$ws.Range("L17").Value = 1494.83724
$ws.Range("N17").Value = -1830.83724
$ws.Range("H19").Value = 1402.875
$ws.Range("I19").Value = 1136.75
$ws.Range("J19").Value = 1491.5834
$ws.Range("K19").Value = 1136.75
$ws.Range("L19").Value = 1491.5834
$ws.Range("M19").Value = -961.75
$ws.Range("N19").Value = -1841.5834
$ws.Range("H29").Value = 6633.3335
$ws.Range("J29").Value = 9500
$ws.Range("L29").Value = 28500
$ws.Range("N29").Value = -29062
$ws.Range("H135").Value = 20940
$ws.Range("I135").Value = 850
$ws.Range("J135").Value = 34333.332
$ws.Range("K135").Value = 7650
$ws.Range("L135").Value = 308999.988
$ws.Range("M135").Value = -5115
$ws.Range("N135").Value = -314069.988
$ws.Range("H137").Value = 4800.9
$ws.Range("I137").Value = 1727.3334
$ws.Range("J137").Value = 9411.25
$ws.Range("K137").Value = 5182.0002
$ws.Range("L137").Value = 28233.75
$ws.Range("M137").Value = -2632.0002
$ws.Range("N137").Value = -33333.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 30632.375
$ws.Range("I31").Value = 3511.5
$ws.Range("K31").Value = 3511.5
$ws.Range("M31").Value = -3217.5
$ws.Range("H61").Value = 31321458
$ws.Range("I61").Value = 45461716
$ws.Range("K61").Value = 45461716
$ws.Range("M61").Value = -45461504
$ws.Range("H63").Value = 2937.889
$ws.Range("I63").Value = 2937.889
$ws.Range("K63").Value = 2937.889
$ws.Range("M63").Value = -2251.889
$ws.Range("H66").Value = 2937.889
$ws.Range("I66").Value = 2937.889
$ws.Range("K66").Value = 14689.445
$ws.Range("M66").Value = -11257.445
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 20000
$ws.Range("K75").Value = 20000
$ws.Range("M75").Value = -19126
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 20000
$ws.Range("K78").Value = 60000
$ws.Range("M78").Value = -55632
$ws.Range("H82").Value = 48537.5
$ws.Range("J82").Value = 48537.5
$ws.Range("L82").Value = 48537.5
$ws.Range("N82").Value = -49259.5
$ws.Range("H85").Value = 48537.5
$ws.Range("J85").Value = 48537.5
$ws.Range("L85").Value = 48537.5
$ws.Range("N85").Value = -51033.5
$ws.Range("H102").Value = 11134.077
$ws.Range("I102").Value = 11674.4
$ws.Range("J102").Value = 9333
$ws.Range("K102").Value = 11674.4
$ws.Range("L102").Value = 9333
$ws.Range("M102").Value = -10052.4
$ws.Range("N102").Value = -12577
$ws.Range("H113").Value = 110000
$ws.Range("J113").Value = 110000
$ws.Range("L113").Value = 110000
$ws.Range("N113").Value = -118678
$ws.Range("H136").Value = 31321458
$ws.Range("I136").Value = 45461716
$ws.Range("K136").Value = 136385148
$ws.Range("M136").Value = -136382598
$ws.Range("H137").Value = 49999
$ws.Range("J137").Value = 49999
$ws.Range("L137").Value = 49999
$ws.Range("N137").Value = -60199

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 46990
$ws.Range("J28").Value = 46990
$ws.Range("L28").Value = 46990
$ws.Range("N28").Value = -47578
$ws.Range("H96").Value = 32201.2
$ws.Range("I96").Value = 12733
$ws.Range("J96").Value = 61403.5
$ws.Range("K96").Value = 12733
$ws.Range("L96").Value = 61403.5
$ws.Range("M96").Value = -9987
$ws.Range("N96").Value = -66895.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2017
$ws.Range("I10").Value = 2017
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2017
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -1878
$ws.Range("H103").Value = 41867.75
$ws.Range("J103").Value = 58737.25
$ws.Range("L103").Value = 58737.25
$ws.Range("N103").Value = -61081.25
$ws.Range("H105").Value = 3282.4614
$ws.Range("I105").Value = 3045.375
$ws.Range("J105").Value = 3661.8
$ws.Range("K105").Value = 3045.375
$ws.Range("L105").Value = 3661.8
$ws.Range("M105").Value = -1298.375
$ws.Range("N105").Value = -7155.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 233.9
$ws.Range("I34").Value = 233.9
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 701.7
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -617.7
$ws.Range("H36").Value = 274.5
$ws.Range("I36").Value = 274.5
$ws.Range("K36").Value = 823.5
$ws.Range("M36").Value = -654.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 138.4
$ws.Range("J2").Value = 145.44444
$ws.Range("L2").Value = 145.44444
$ws.Range("N2").Value = -371.44444
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5
$ws.Range("K19").Value = 5
$ws.Range("M19").Value = 283
$ws.Range("H62").Value = 99954.664
$ws.Range("J62").Value = 99954.664
$ws.Range("L62").Value = 99954.664
$ws.Range("N62").Value = -101326.664
$ws.Range("H65").Value = 99954.664
$ws.Range("J65").Value = 99954.664
$ws.Range("L65").Value = 299863.992
$ws.Range("N65").Value = -306727.992

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 116000
$ws.Range("J81").Value = 116000
$ws.Range("L81").Value = 116000
$ws.Range("N81").Value = -117996
$ws.Range("H84").Value = 116000
$ws.Range("J84").Value = 116000
$ws.Range("L84").Value = 348000
$ws.Range("N84").Value = -357984
$ws.Range("H95").Value = 27449.5
$ws.Range("J95").Value = 27449.5
$ws.Range("L95").Value = 27449.5
$ws.Range("N95").Value = -32941.5
$ws.Range("H132").Value = 155509.72
$ws.Range("I132").Value = 113014
$ws.Range("J132").Value = 232002
$ws.Range("K132").Value = 339042
$ws.Range("L132").Value = 696006
$ws.Range("M132").Value = -336512
$ws.Range("N132").Value = -701066

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4142.5713
$ws.Range("J96").Value = 6133
$ws.Range("L96").Value = 6133
$ws.Range("N96").Value = -8879

